$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate rows 2-5 with the updated TPM-derived NATMI ligand-receptor values.
# Row 2-3: Sending cluster = ECs (was previously MuSCs-only in rows 2-3)
# Row 4-5: Sending cluster = MuSCs (new rows, Target cluster = ECs / FAPs)

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nppc"
$ws.Range("C2").Value = "Npr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.224617
$ws.Range("H2").Value = 0.673851
$ws.Range("I2").Value = 0.7269197784238318
$ws.Range("J2").Value = 0.7269197784238318
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2689956666666667
$ws.Range("N2").Value = 0.806987
$ws.Range("O2").Value = 0.1265890356442505
$ws.Range("P2").Value = 0.1265890356442505
$ws.Range("Q2").Value = 0.06042099965966666
$ws.Range("R2").Value = 0.543788996937
$ws.Range("S2").Value = 0.09202007374140514
$ws.Range("T2").Value = 0.09202007374140514

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nppc"
$ws.Range("C3").Value = "Npr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.224617
$ws.Range("H3").Value = 0.673851
$ws.Range("I3").Value = 0.7269197784238318
$ws.Range("J3").Value = 0.7269197784238318
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.855956666666666
$ws.Range("N3").Value = 5.567869999999999
$ws.Range("O3").Value = 0.8734109643557494
$ws.Range("P3").Value = 0.8734109643557494
$ws.Range("Q3").Value = 0.4168794185966666
$ws.Range("R3").Value = 3.751914767369999
$ws.Range("S3").Value = 0.6348997046824266
$ws.Range("T3").Value = 0.6348997046824266

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Nppc"
$ws.Range("C4").Value = "Npr3"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.08438133333333332
$ws.Range("H4").Value = 0.253144
$ws.Range("I4").Value = 0.2730802215761681
$ws.Range("J4").Value = 0.2730802215761681
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2689956666666667
$ws.Range("N4").Value = 0.806987
$ws.Range("O4").Value = 0.1265890356442505
$ws.Range("P4").Value = 0.1265890356442505
$ws.Range("Q4").Value = 0.02269821301422222
$ws.Range("R4").Value = 0.204283917128
$ws.Range("S4").Value = 0.03456896190284538
$ws.Range("T4").Value = 0.03456896190284538

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Nppc"
$ws.Range("C5").Value = "Npr3"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.08438133333333332
$ws.Range("H5").Value = 0.253144
$ws.Range("I5").Value = 0.2730802215761681
$ws.Range("J5").Value = 0.2730802215761681
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.855956666666666
$ws.Range("N5").Value = 5.567869999999999
$ws.Range("O5").Value = 0.8734109643557494
$ws.Range("P5").Value = 0.8734109643557494
$ws.Range("Q5").Value = 0.1566080981422222
$ws.Range("R5").Value = 1.40947288328
$ws.Range("S5").Value = 0.2385112596733227
$ws.Range("T5").Value = 0.2385112596733227
